$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Portal Check Added and Asana Updates"
# The SKU / Vendor-Part / Units-Received / Units-Invoiced / PO-Unit-Cost /
# Invoice-Unit-Cost / Extended-Cost-Variance data block (columns I,J,L,M,N,O,P)
# rotates across the three data rows:
#   new row 2 <- old row 3
#   new row 3 <- old row 4
#   new row 4 <- old row 2
# Column K ("SKU DESCRIPTION", a single space) is identical on every row, so it
# is left untouched. Copy/paste (rather than re-typing values) keeps the cells
# as shared-string text instead of turning them into numbers.

# Stash old row 2's data in a scratch range far outside the used area.
$ws.Range("I2:P2").Copy($ws.Range("AA2"))

# Shift row 3 -> row 2, row 4 -> row 3.
$ws.Range("I3:P3").Copy($ws.Range("I2"))
$ws.Range("I4:P4").Copy($ws.Range("I3"))

# Old row 2's data becomes the new row 4.
$ws.Range("AA2:AH2").Copy($ws.Range("I4"))

# Clean up the scratch range.
$ws.Range("AA2:AH2").ClearContents()
